$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A77").Value = "mail_new_course_registration_requested_body"
$ws.Range("B77").Value = "{0}, Hemos recibido tú solicitud de registro del curso {1}, en breve completaremos el registro!"
$ws.Range("C77").Value = "{0}, We have received your registration request for the course {1}, we will complete the registration shortly!"

$ws.Range("A78").Value = "mail_certificate_course_disabled_body"
$ws.Range("A79").Value = "mail_certificate_course_enabled_body"

$ws.Range("B78").Value = "{0}, El curso {1} ha sido deshabilitado en TCS, deberás volver a habilitarlo para que nuevos certificados puedan ser emitidos."
$ws.Range("B79").Value = "{0}, El curso {1} ha sido habilitado en TCS, nuevos certificados podrán ser emitidos."

$ws.Range("C78").Value = "{0}, The course {1} has been disabled in TCS, you must re-enable it so that new certificates can be issued."
$ws.Range("C79").Value = "{0}, The course {1} has been enabled in TCS, new certificates may be issued."

$ws.Range("C77:C79").Font.Size = 11

$ws.Range("C79").Select()
